$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.0864417552948
$ws.Range("B1").Value = 1.231530427932739
$ws.Range("C1").Value = 1.181447863578796
$ws.Range("D1").Value = 1.405853629112244
$ws.Range("E1").Value = 1.274385333061218
